$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.042.11'
$ws.Range("E2").Value = '  -6.60%  '
$ws.Range("D3").Value = '2.551.81'
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '298.38'
$ws.Range("E5").Value = '  -4.41%  '
$ws.Range("D6").Value = '92.12'
$ws.Range("E6").Value = '  -7.33%  '
$ws.Range("E7").Value = '  -3.84%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.552'
$ws.Range("E9").Value = '  -5.44%  '
$ws.Range("D10").Value = '36.06'
$ws.Range("E10").Value = '  -7.82%  '
$ws.Range("D11").Value = '0.0807'
$ws.Range("E11").Value = '  -4.36%  '
$ws.Range("D12").Value = '7.74'
$ws.Range("E12").Value = '  -5.32%  '
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '2.936.77'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = '2.538.88'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").Value = '0.872'
$ws.Range("E16").Value = '  -4.97%  '
$ws.Range("D17").Value = '14.18'
$ws.Range("E17").Value = '  -4.69%  '
$ws.Range("D18").Value = '42.973.39'
$ws.Range("E18").Value = '  -7.01%  '
$ws.Range("D19").Value = '6.67'
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("D20").Value = '0.0₃0978'
$ws.Range("E20").Value = '  -4.08%  '
$ws.Range("D21").Value = '12.54'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").Value = '72.17'
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").Value = '261.23'
$ws.Range("E23").Value = '  -9.80%  '
$ws.Range("D24").Value = '2.92'
$ws.Range("E24").Value = '  -4.95%  '
$ws.Range("D25").Value = '29.68'
$ws.Range("E25").Value = '  +1.55%  '
$ws.Range("D26").Value = '2.14'
$ws.Range("E26").Value = '  -4.06%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '10.09'
$ws.Range("D29").Value = '2.13'
$ws.Range("E29").Value = '  -3.73%  '
$ws.Range("D30").Value = '''36.60'
$ws.Range("E30").Value = '  -6.79%  '
$ws.Range("D31").Value = '''6.00'
$ws.Range("E31").Value = '  -4.18%  '
$ws.Range("D32").Value = '153.31'
$ws.Range("E32").Value = '  -1.95%  '
$ws.Range("D33").Value = '2.18'
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("D34").Value = '''3.40'
$ws.Range("E34").Value = '  -5.23%  '
$ws.Range("D35").Value = '2.72'
$ws.Range("E35").Value = '  -2.37%  '
$ws.Range("D36").Value = '0.0794'
$ws.Range("E36").Value = '  -5.53%  '
$ws.Range("D37").Value = '0.114'
$ws.Range("E37").Value = '  -6.78%  '
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("D39").Value = '24.09'
$ws.Range("E39").Value = '  +14.91%  '
$ws.Range("D40").Value = '16.82'
$ws.Range("E40").Value = '  +7.07%  '
$ws.Range("D41").Value = '3.47'
$ws.Range("E41").Value = '  -3.46%  '
$ws.Range("D42").Value = '0.0312'
$ws.Range("E42").Value = '  -6.30%  '
$ws.Range("D43").Value = '3.84'
$ws.Range("E43").Value = '  -4.31%  '
$ws.Range("D44").Value = '2.077.73'
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("D45").Value = '0.997'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '85.76'
$ws.Range("E46").Value = '  -11.72%  '
$ws.Range("E47").Value = '  +3.30%  '
$ws.Range("D48").Value = '2.794.13'
$ws.Range("D49").Value = '1.71'
$ws.Range("E49").Value = '  -2.12%  '
$ws.Range("D50").Value = '''104.40'
$ws.Range("E50").Value = '  -4.10%  '
$ws.Range("E51").Value = '  -8.60%  '
